# Update: learning guide 02 resources
# Ventas-mayo-2024 -> renamed from "figuras" (LEGO minifig colors) theme
# to "set Lego" (LEGO set names) theme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title ---
$ws.Range("B1").Value = "Venta de set Lego mayo 2024"

# --- Column header (B3): "Figura" -> "Set" ---
$ws.Range("B3").Value = "Set"

# --- Product names in column B (Lego colors -> Lego set names) ---
$ws.Range("B6").Value  = "BASH!"
$ws.Range("B7").Value  = "Racing Yacht"
$ws.Range("B8").Value  = "First Responder"
$ws.Range("B9").Value  = "Hovercraft"
$ws.Range("B12").Value = "Rally Car"
$ws.Range("B13").Value = "Mack Anthem"
$ws.Range("B14").Value = "First Responder"
$ws.Range("B15").Value = "Forest Machine"
$ws.Range("B17").Value = "Rough Terrain Crane"
$ws.Range("B18").Value = "Bugatti Chiron"
$ws.Range("B19").Value = "Hook Loader"

# --- Give the title its own (explicit) font style ---
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 11

# --- Widen column B to fit the longer set names ---
$ws.Columns("B").ColumnWidth = 20.6

# --- Stray formatted (empty) cell that shows up past the table, with an
#     underlined font, matching the extended used range (K12) ---
$ws.Range("K12").Font.Underline = $true

# --- Selection moved one column/row further right, following the widened
#     used range ---
$ws.Range("L14").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
